# Append " *" to the header labels to mark them as required fields,
# and move the active selection to D2 (from D10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Fund *"
$ws.Range("B1").Value = "Name *"
$ws.Range("C1").Value = "Percentage Called *"
$ws.Range("D1").Value = "Due Date *"

$ws.Range("D2").Select()
